# Cigna flow - add URL_Cigna configuration row to the Settings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row above the existing "URL_Aetna" row (row 3), shifting the
# remaining configuration rows (and the blank rows / B14 marker below them)
# down by one.
$ws.Rows.Item(3).Insert()

# Populate the new row with the Cigna website configuration entry.
# Set Name (A) then Description (C) then Value (B) so the new shared
# strings are appended to the shared string table in that order.
$ws.Range("A3").Value = "URL_Cigna"
$ws.Range("C3").Value = "Cigna Website URL"
$ws.Range("B3").Value = "https://hcpdirectory.cigna.com/web/public/consumer/directory/search?consumerCode=HDC001"

# Reflect the author's new active cell selection.
$ws.Range("B7").Select() | Out-Null
